$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns K/L (ArrangementOffsetX / ArrangementOffsetY) with their
# type-row ("float") metadata, mirroring the existing table layout.
$ws.Range("K3").Value = "ArrangementOffsetX"
$ws.Range("L3").Value = "ArrangementOffsetY"
$ws.Range("K4").Value = "float"
$ws.Range("L4").Value = "float"

# Data rows (5-15) use vertically-centered cells like the rest of the table.
$xlVAlignCenter = -4108
$ws.Range("K5:L15").VerticalAlignment = $xlVAlignCenter

# Update the SlotCount (F) values and D13's display name, then populate the
# new ArrangementOffsetX/Y columns for every data row with their defaults.
$ws.Range("F5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0

$ws.Range("F6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0

$ws.Range("F7").Value = 1
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

$ws.Range("F8").Value = 3
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0

$ws.Range("F9").Value = 1
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0

$ws.Range("F10").Value = 3
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

$ws.Range("F11").Value = 3
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0

$ws.Range("F12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0

$ws.Range("D13").Value = "Fresh Display 2"
$ws.Range("F13").Value = 1
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0

$ws.Range("F14").Value = 2
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0

$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0

# Restore the selection left behind by the editing session.
$ws.Range("E17").Select() | Out-Null
